$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Grayson Allen",      "PG,SG,SF", "Phoenix Suns"),
    @("Michael Porter Jr.", "SF,PF",    "Denver Nuggets"),
    @("Jaden McDaniels",    "SF,PF",    "Minnesota Timberwolves"),
    @("Toumani Camara",     "SG,SF,PF", "Portland Trail Blazers"),
    @("Domantas Sabonis",   "C",        "Sacramento Kings"),
    @("Alperen Sengün",     "C",        "Houston Rockets"),
    @("Kristaps Porzingis", "PF,C",     "Boston Celtics"),
    @("Josh Hart",          "SG,SF,PF", "New York Knicks"),
    @("Kelly Oubre Jr.",    "SG,SF",    "Philadelphia 76ers"),
    @("Aaron Wiggins",      "SG,SF",    "Oklahoma City Thunder"),
    @("Dyson Daniels",      "PG,SG,SF", "Atlanta Hawks"),
    @("Malik Beasley",      "SG,SF",    "Detroit Pistons"),
    @("Devin Vassell",      "SG,SF",    "San Antonio Spurs"),
    @("Donovan Mitchell",   "PG,SG",    "Cleveland Cavaliers"),
    @("De'Andre Hunter",    "SF,PF",    "Cleveland Cavaliers"),
    @("Julius Randle",      "PF,C",     "Minnesota Timberwolves"),
    @("Cam Thomas",         "SG,SF",    "Brooklyn Nets"),
    @("Guerschon Yabusele", "PF,C",     "Philadelphia 76ers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
